# Update profit.py after running on 2025-08-30
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: append the new day's row (row 13) ---
# Write the new date as text (matching the existing plain-text date cells
# above it) instead of letting Excel auto-convert the string into a date
# serial number. Force Text format first, write the value, then restore
# the default (unstyled) cell formatting by pasting it in from a sibling
# date cell so no extra style entries get introduced into the workbook.
$ws1.Range("A13").NumberFormat = "@"
$ws1.Range("A13").Value = "08/30/2025"
$ws1.Range("A12").Copy()
$ws1.Range("A13").PasteSpecial(-4122)

$ws1.Range("B13").Value = 11447.83

# --- Sheet2: overwrite the single summary row with the new date/ratios ---
$ws2.Range("A1").NumberFormat = "@"
$ws2.Range("A1").Value = "08/30/2025"
$ws1.Range("A12").Copy()
$ws2.Range("A1").PasteSpecial(-4122)

$ws2.Range("B1").Value = 0.1105078176324462
$ws2.Range("C1").Value = 0.8894921823675538
